$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Student test 1"
$ws.Range("D3").Value = "18 Maret 2025 ya kakk, jangan sampai lupa :)"
$ws.Range("C3").Value = "Python Intermediate"
$ws.Range("D2").Value = "bisa dengan menggunakan print('halo') kak"
$ws.Range("B3").Value = "Student test 1"

$ws.Range("D2").Select()
